$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 35714588
$ws.Range("I2").Value = 35714588
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 35714588
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -35714475
$ws.Range("N2").ClearContents()
$ws.Range("H53").Value = 92093.91
$ws.Range("I53").Value = 144324.72
$ws.Range("J53").Value = 690
$ws.Range("K53").Value = 144324.72
$ws.Range("L53").Value = 690
$ws.Range("M53").Value = -143687.72
$ws.Range("N53").Value = -1964
$ws.Range("H92").Value = 43403460
$ws.Range("I92").Value = 2315413.8
$ws.Range("K92").Value = 2315413.8
$ws.Range("M92").Value = -2314165.8
$ws.Range("H107").Value = 7813017
$ws.Range("I107").Value = 10869854
$ws.Range("J107").Value = 1100.8889
$ws.Range("K107").Value = 10869854
$ws.Range("L107").Value = 1100.8889
$ws.Range("M107").Value = -10867934
$ws.Range("N107").Value = -4940.8889
$ws.Range("H112").Value = 55574828
$ws.Range("I112").Value = 650
$ws.Range("J112").Value = 125042550
$ws.Range("K112").Value = 1950
$ws.Range("L112").Value = 375127650
$ws.Range("M112").Value = -842
$ws.Range("N112").Value = -375129866
$ws.Range("H129").Value = 894.202
$ws.Range("I129").Value = 460.875
$ws.Range("J129").Value = 977.7349
$ws.Range("K129").Value = 1382.625
$ws.Range("L129").Value = 2933.2047
$ws.Range("M129").Value = 3617.375
$ws.Range("N129").Value = -12933.2047
$ws.Range("H137").Value = 1321.05
$ws.Range("I137").Value = 1236.2142
$ws.Range("J137").Value = 1519
$ws.Range("K137").Value = 3708.6426
$ws.Range("L137").Value = 4557
$ws.Range("M137").Value = -1158.6426
$ws.Range("N137").Value = -9657

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H17").Value = 4490
$ws.Range("J17").Value = 4490
$ws.Range("L17").Value = 4490
$ws.Range("N17").Value = -4836
$ws.Range("H32").Value = 5485.222
$ws.Range("I32").Value = 3774.1714
$ws.Range("J32").Value = 16373.728
$ws.Range("K32").Value = 3774.1714
$ws.Range("L32").Value = 16373.728
$ws.Range("M32").Value = -3487.1714
$ws.Range("N32").Value = -16947.728
$ws.Range("H41").Value = 3685.3333
$ws.Range("I41").Value = 3685.3333
$ws.Range("K41").Value = 3685.3333
$ws.Range("M41").Value = -3271.3333

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 100001250
$ws.Range("I99").Value = 125001180
$ws.Range("J99").Value = 1500
$ws.Range("K99").Value = 125001180
$ws.Range("L99").Value = 1500
$ws.Range("M99").Value = -124999682
$ws.Range("N99").Value = -4496
$ws.Range("H125").Value = 58113.332
$ws.Range("J125").Value = 58113.332
$ws.Range("L125").Value = 58113.332
$ws.Range("N125").Value = -67953.33199999999

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1522.2941
$ws.Range("I16").Value = 1089.8572
$ws.Range("J16").Value = 1825
$ws.Range("K16").Value = 1089.8572
$ws.Range("L16").Value = 1825
$ws.Range("M16").Value = -802.8571999999999
$ws.Range("N16").Value = -2399
$ws.Range("H31").Value = 13516976
$ws.Range("I31").Value = 1631.2858
$ws.Range("J31").Value = 31255866
$ws.Range("K31").Value = 1631.2858
$ws.Range("L31").Value = 31255866
$ws.Range("M31").Value = -1336.2858
$ws.Range("N31").Value = -31256456
$ws.Range("H34").Value = 13516976
$ws.Range("I34").Value = 1631.2858
$ws.Range("J34").Value = 31255866
$ws.Range("K34").Value = 1631.2858
$ws.Range("L34").Value = 31255866
$ws.Range("M34").Value = -1429.2858
$ws.Range("N34").Value = -31256270
$ws.Range("H113").Value = 1522.2941
$ws.Range("I113").Value = 1089.8572
$ws.Range("J113").Value = 1825
$ws.Range("K113").Value = 1089.8572
$ws.Range("L113").Value = 1825
$ws.Range("M113").Value = 1080.1428
$ws.Range("N113").Value = -6165
$ws.Range("H132").Value = 8001968.5
$ws.Range("I132").Value = 9092632
$ws.Range("K132").Value = 27277896
$ws.Range("M132").Value = -27275366

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 3448405
$ws.Range("I12").Value = 11111226
$ws.Range("J12").Value = 135.55
$ws.Range("K12").Value = 33333678
$ws.Range("L12").Value = 406.65
$ws.Range("M12").Value = -33333505
$ws.Range("N12").Value = -752.6500000000001
$ws.Range("H18").Value = 316.66666
$ws.Range("I18").Value = 445
$ws.Range("J18").Value = 60
$ws.Range("K18").Value = 1335
$ws.Range("L18").Value = 180
$ws.Range("M18").Value = -1166
$ws.Range("N18").Value = -518
$ws.Range("H23").Value = 6250094
$ws.Range("J23").Value = 97.90909000000001
$ws.Range("L23").Value = 293.72727
$ws.Range("N23").Value = -763.7272700000001
$ws.Range("H29").Value = 1355
$ws.Range("I29").Value = 30
$ws.Range("J29").Value = 1796.6666
$ws.Range("K29").Value = 90
$ws.Range("L29").Value = 5389.9998
$ws.Range("M29").Value = 187
$ws.Range("N29").Value = -5943.9998
$ws.Range("H31").Value = 500
$ws.Range("J31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("N31").ClearContents()
$ws.Range("H38").Value = 7692344.5
$ws.Range("I38").Value = 12500022
$ws.Range("J38").Value = 60.4
$ws.Range("K38").Value = 37500066
$ws.Range("L38").Value = 181.2
$ws.Range("M38").Value = -37499719
$ws.Range("N38").Value = -875.2
$ws.Range("H80").Value = 2000
$ws.Range("I80").Value = 500
$ws.Range("J80").Value = 2750
$ws.Range("K80").Value = 1500
$ws.Range("L80").Value = 8250
$ws.Range("M80").Value = -564
$ws.Range("N80").Value = -10122
$ws.Range("H83").Value = 2000
$ws.Range("I83").Value = 500
$ws.Range("J83").Value = 2750
$ws.Range("K83").Value = 4500
$ws.Range("L83").Value = 24750
$ws.Range("M83").Value = 180
$ws.Range("N83").Value = -34110
$ws.Range("H87").Value = 1600
$ws.Range("I87").Value = 1600
$ws.Range("K87").Value = 4800
$ws.Range("M87").Value = -3552
$ws.Range("H90").Value = 1600
$ws.Range("I90").Value = 1600
$ws.Range("K90").Value = 14400
$ws.Range("M90").Value = -8160
$ws.Range("H107").Value = 521.0526
$ws.Range("I107").Value = 337.2
$ws.Range("J107").Value = 586.7143
$ws.Range("K107").Value = 1011.6
$ws.Range("L107").Value = 1760.1429
$ws.Range("M107").Value = 908.4000000000001
$ws.Range("N107").Value = -5600.1429
$ws.Range("H113").Value = 4615868.5
$ws.Range("I113").Value = 6250474.5
$ws.Range("J113").Value = 2000498.8
$ws.Range("K113").Value = 18751423.5
$ws.Range("L113").Value = 6001496.4
$ws.Range("M113").Value = -18749253.5
$ws.Range("N113").Value = -6005836.4
$ws.Range("H132").Value = 13890315
$ws.Range("I132").Value = 1175
$ws.Range("J132").Value = 27779454
$ws.Range("K132").Value = 10575
$ws.Range("L132").Value = 250015086
$ws.Range("M132").Value = -8045
$ws.Range("N132").Value = -250020146

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 20140.334
$ws.Range("I15").Value = 20140.334
$ws.Range("K15").Value = 20140.334
$ws.Range("M15").Value = -19852.334
$ws.Range("H81").Value = 20140.334
$ws.Range("I81").Value = 20140.334
$ws.Range("K81").Value = 20140.334
$ws.Range("M81").Value = -19142.334
$ws.Range("H84").Value = 20140.334
$ws.Range("I84").Value = 20140.334
$ws.Range("K84").Value = 60421.00199999999
$ws.Range("M84").Value = -55429.00199999999
$ws.Range("H132").Value = 9806754
$ws.Range("I132").Value = 11907201
$ws.Range("J132").Value = 4666
$ws.Range("K132").Value = 35721603
$ws.Range("L132").Value = 13998
$ws.Range("M132").Value = -35719073
$ws.Range("N132").Value = -19058

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3394.95
$ws.Range("I61").Value = 2589.4
$ws.Range("J61").Value = 4200.5
$ws.Range("K61").Value = 2589.4
$ws.Range("L61").Value = 4200.5
$ws.Range("M61").Value = -2387.4
$ws.Range("N61").Value = -4604.5
$ws.Range("H113").Value = 3394.95
$ws.Range("I113").Value = 2589.4
$ws.Range("J113").Value = 4200.5
$ws.Range("K113").Value = 2589.4
$ws.Range("L113").Value = 4200.5
$ws.Range("M113").Value = -419.4000000000001
$ws.Range("N113").Value = -8540.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H57").Value = 49600
$ws.Range("I57").Value = 49200
$ws.Range("J57").Value = 50000
$ws.Range("K57").Value = 49200
$ws.Range("L57").Value = 50000
$ws.Range("M57").Value = -48446
$ws.Range("N57").Value = -51508
$ws.Range("H68").Value = 42300
$ws.Range("J68").Value = 42300
$ws.Range("L68").Value = 42300
$ws.Range("N68").Value = -43922
$ws.Range("H71").Value = 42300
$ws.Range("J71").Value = 42300
$ws.Range("L71").Value = 126900
$ws.Range("N71").Value = -135012
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()
$ws.Range("H80").Value = 2000
$ws.Range("I80").Value = 2000
$ws.Range("K80").Value = 2000
$ws.Range("M80").Value = -1002
$ws.Range("H83").Value = 2000
$ws.Range("I83").Value = 2000
$ws.Range("K83").Value = 6000
$ws.Range("M83").Value = -1008
$ws.Range("H136").Value = 3547005.2
$ws.Range("I136").Value = 504.9091
$ws.Range("J136").Value = 11906613
$ws.Range("K136").Value = 1514.7273
$ws.Range("L136").Value = 35719839
$ws.Range("M136").Value = 1035.2727
$ws.Range("N136").Value = -35724939
